$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Mean"
$ws.Range("B5").Value = 0.03172419637273724
$ws.Range("C5").Value = 0.01317945728601877
$ws.Range("D5").Value = 0.04619917179409822
$ws.Range("E5").Value = 0.0399047008522001
$ws.Range("F5").Value = 0.07814089129967296
$ws.Range("G5").Value = 0.1238611781735173
